# Apply the changes described by the diff:
# - model sheet: is_distributed / is_override fields change type from "string" to "boolean"
# - choices sheet: true_false choice values change from text "true"/"false" to numeric 1/0
# - active sheet/selection changes: "model" becomes the active tab (was "choices"),
#   with the selection on "choices" moving to G19 and on "model" moving to A14.

$wb = $excel.ActiveWorkbook

# --- model sheet: update field types to "boolean" ---
$model = $wb.Worksheets.Item("model")
$model.Range("A11").Value = "boolean"   # is_distributed
$model.Range("A12").Value = "boolean"   # is_override

# --- choices sheet: true_false values become numeric 1/0 ---
$choices = $wb.Worksheets.Item("choices")
$choices.Range("B2").Value = 1          # was "true"
$choices.Range("B3").Value = 0          # was "false"

# --- update view/selection state ---
# Select on "choices" sheet first (it loses the active/tabSelected flag).
$choices.Activate()
$choices.Range("G19").Select()

# Then activate "model" sheet last so it becomes the active tab, with its own selection.
$model.Activate()
$model.Range("A14").Select()
